$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.112.51'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.799.50'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '222.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.286'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('E10').Value = '  +4.31%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.056.74'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.787.67'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.72'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.125.77'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.28'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.10'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0786'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.86'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.56'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0526'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.51'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.415.70'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.649'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.51%  '
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.945'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.27'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('E43').Value = '  +4.57%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.50'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.06%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0494'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.955.70'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.91'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0124'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.14%  '
